$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value as genuine TEXT (matches the
# workbook's existing convention of caching amount columns as text, e.g.
# "22016.00") without leaving a stray quote-prefix style on the cell.
function Set-TextValue($sheet, $addr, $val) {
    $rng = $sheet.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

# Row 9 - "Short point" line
$ws.Range("C9").Value = 44
Set-TextValue $ws "G9" "11264.00"

# Row 10 - "Medium point" line
$ws.Range("C10").Value = 85
Set-TextValue $ws "G10" "40120.00"

# Row 11 - "Long point" line
$ws.Range("C11").Value = 25
Set-TextValue $ws "G11" "16550.00"

# Row 12 - header row for plug point item (qty only)
$ws.Range("C12").Value = 13

# Row 13 - "On board" line
$ws.Range("C13").Value = 4
Set-TextValue $ws "G13" "544.00"

# Row 14 - switch item line
$ws.Range("C14").Value = 16
Set-TextValue $ws "G14" "368.00"

# Row 15 - "Total" line (qty only)
$ws.Range("C15").Value = 18

# Row 16 - "Add Tender Premium" line (qty only)
$ws.Range("C16").Value = 51

# Row 17 - "Grand Total" line (qty only)
$ws.Range("C17").Value = 12

# Row 19 - Grand Total Rs. summary
Set-TextValue $ws "G19" "68846.00"
Set-TextValue $ws "H19" "68846.00"

# Row 21 - NET PAYABLE AMOUNT Rs. summary
Set-TextValue $ws "G21" "68846.00"
Set-TextValue $ws "H21" "68846.00"
